# Update folder paths to point at the new repo locations, and add a
# raw_data column pointing at the BIDS sourcedata tree.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldBase = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\OSF_files\EEG_derivatives\sub-"
$rawBase = "I:\SCIENCE-NEXS-neurolab\PROJECTS\PLAYMORE\EEG_project1\Analyses\BIDS_data\bids\sourcedata\sub-"

# New header for column H
$ws.Cells.Item(1, 8).Value = "raw_data"

# Last data row is 22 (rows 2..22 hold subjects sub-01..sub-21, in order)
$lastRow = 22

for ($r = 2; $r -le $lastRow; $r++) {
    $subLabel = $ws.Cells.Item($r, 7).Value()       # e.g. "sub-01"
    $num = $subLabel.Substring(4)                   # "01"

    $ws.Cells.Item($r, 2).Value = "$oldBase$num"
    $ws.Cells.Item($r, 8).Value = "$rawBase$num"
}

$ws.Range("G16").Select()
